$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Cell values (new compound swapped in for row 2, GHS hazard info added,
#    physical-property blobs reformatted with one property per line, and a
#    new "Hazards" entry added at C3).
# ---------------------------------------------------------------------------

$ws.Range("A1").Value = "Compound"
$ws.Range("B1").Value = "Physical Properties"
$ws.Range("C1").Value = "Hazards"
$ws.Range("D1").Value = "First Aids"

$ws.Range("A2").Value = "Hydrochloric acid"
$ws.Range("B2").Value = "ordor: Pungent, irritating odor;`nBoiling Point: 123 °F at 760 mm Hg (USCG, 1999);`nMelting Point: -174.6 °F (Melting point is -13.7° F for a 39.17% weight/weight solution.) (EPA, 1998);`nsolubility: 82.3 g/100 g at 32° F (NTP, 1992);`nDensity: 1.05 at 59 °F for 10.17% weight/weight solution (EPA, 1998);`n"
$ws.Range("C2").Value = "H314: Causes severe skin burns and eye damage [Danger Skin corrosion/irritation]`nH331: Toxic if inhaled [Danger Acute toxicity, inhalation]"
$ws.Range("D2").Value = "inhalation: Fresh air, rest. Half-upright position. Artificial respiration may be needed. Refer immediately for medical attention.`nskin: Wear protective gloves when administering first aid. First rinse with plenty of water for at least 15 minutes, then remove contaminated clothes and rinse again. Refer immediately for medical attention.`neye: Rinse with plenty of water for several minutes (remove contact lenses if easily possible). Refer immediately for medical attention.`n"

$ws.Range("A3").Value = "Calcium carbonate"
$ws.Range("B3").Value = "ordor: Odorless;`nBoiling Point: Decomposes (NIOSH, 2016);`nMelting Point: 1517 to 2442 °F (Decomposes) (NIOSH, 2016);`nsolubility: 0.001 % (NIOSH, 2016);`nDensity: 2.7 to 2.95 (NIOSH, 2016);`n"
$ws.Range("C3").Value = "Not Classified`nReported as not meeting GHS hazard criteria by 3306 of 3614 companies (only ~ 8.5% companies provided GHS information). For more detailed information, please visit ECHA C&L website"
$ws.Range("D3").Value = "inhalation: Fresh air.`nskin: Rinse skin with plenty of water or shower.`neye: Rinse with plenty of water (remove contact lenses if easily possible).`ningestion: Rinse mouth.`n"

$ws.Range("A4").Value = "Sodium chloride"
$ws.Range("B4").Value = "Boiling Point: 2575 °F at 760 mm Hg (NTP, 1992);`nMelting Point: 1474 °F (NTP, 1992);`nsolubility: greater than or equal to 100 mg/mL at 68° F (NTP, 1992);`nDensity: 2.165 at 77 °F (NTP, 1992);`n"
$ws.Range("C4").Value = "Aggregated GHS information provided by 2341 companies from 13 notifications to the ECHA C&L Inventory.`nReported as not meeting GHS hazard criteria by 1699 of 2341 companies. For more detailed information, please visit ECHA C&L website`nOf the 9 notification(s) provided by 642 of 2341 companies with hazard statement code(s):`nH319 (100%): Causes serious eye irritation [Warning Serious eye damage/eye irritation]`nInformation may vary between notifications depending on impurities, additives, and other factors. The percentage value in parenthesis indicates the notified classification ratio from companies that provide hazard codes. Only hazard codes with percentage values above 10% are shown."
$ws.Range("D4").Value = "EYES: First check the victim for contact lenses and remove if present. Flush victim's eyes with water or normal saline solution for 20 to 30 minutes while simultaneously calling a hospital or poison control center. Do not put any ointments, oils, or medication in the victim's eyes without specific instructions from a physician. IMMEDIATELY transport the victim after flushing eyes to a hospital even if no symptoms (such as redness or irritation) develop. SKIN: IMMEDIATELY flood affected skin with water while removing and isolating all contaminated clothing. Gently wash all affected skin areas thoroughly with soap and water. If symptoms such as redness or irritation develop, IMMEDIATELY call a physician and be prepared to transport the victim to a hospital for treatment. INHALATION: IMMEDIATELY leave the contaminated area; take deep breaths of fresh air. If symptoms (such as wheezing, coughing, shortness of breath, or burning in the mouth, throat, or chest) develop, call a physician and be prepared to transport the victim to a hospital. Provide proper respiratory protection to rescuers entering an unknown atmosphere. Whenever possible, Self-Contained Breathing Apparatus (SCBA) should be used; if not available, use a level of protection greater than or equal to that advised under Protective Clothing. INGESTION: DO NOT INDUCE VOMITING. If the victim is conscious and not convulsing, give 1 or 2 glasses of water to dilute the chemical and IMMEDIATELY call a hospital or poison control center. Be prepared to transport the victim to a hospital if advised by a physician. If the victim is convulsing or unconscious, do not give anything by mouth, ensure that the victim's airway is open and lay the victim on his/her side with the head lower than the body. DO NOT INDUCE VOMITING. IMMEDIATELY transport the victim to a hospital. (NTP, 1992)"

# ---------------------------------------------------------------------------
# 2. Formatting: bold/centered wrapped header row, centered wrapped compound
#    column, left-aligned wrapped data columns.
# ---------------------------------------------------------------------------

$header = $ws.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108
$header.WrapText = $true

$colA = $ws.Range("A2:A4")
$colA.HorizontalAlignment = -4108
$colA.VerticalAlignment = -4108
$colA.WrapText = $true

$dataCols = $ws.Range("B2:D4")
$dataCols.HorizontalAlignment = -4131
$dataCols.VerticalAlignment = -4108
$dataCols.WrapText = $true

# ---------------------------------------------------------------------------
# 3. Row heights / column widths to fit the longer text blocks.
# ---------------------------------------------------------------------------

$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 120
$ws.Rows.Item(3).RowHeight = 105
$ws.Rows.Item(4).RowHeight = 370

$ws.Columns.Item(1).ColumnWidth = 17.1666666666667
$ws.Columns.Item(2).ColumnWidth = 68.6666666666667
$ws.Columns.Item(3).ColumnWidth = 63.6666666666667
$ws.Columns.Item(4).ColumnWidth = 80.3307291666667

# ---------------------------------------------------------------------------
# 4. Font face (workbook-wide font switched to a CJK-capable family).
# ---------------------------------------------------------------------------

$ws.Cells.Font.Name = "宋体"

[void]$ws.Range("B4").Select()
